$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update tarifas_clp (E) and hh_totales... columns with new values
$ws.Range("E2").Value = 49176755
$ws.Range("F2").Value = 6

$ws.Range("E3").Value = 29929500
$ws.Range("F3").Value = 6

$ws.Range("E4").Value = 15962400
$ws.Range("F4").Value = 6

# Update the active selection to match the saved view state
$ws.Range("G6").Select()

$wb.Save()
